$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: title -> new text, merged A1:I1, bigger row height, centered/wrapped
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Shuakhevi Municipality"
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------------------
# Row 2: unchanged text, but row height becomes the default (no custom height)
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# Row 4: relabel + fill in real numbers (was placeholder text row)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 968
$ws.Range("C4").Value = 958
$ws.Range("D4").Value = 940
$ws.Range("E4").Value = 972
$ws.Range("F4").Value = 980
$ws.Range("G4").Value = 969
$ws.Range("H4").Value = 964
$ws.Range("I4").Value = 945
$ws.Range("B4:I4").NumberFormat = "#\ ##0"
$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------------------
# Row 5: was the merged "Number of disability persons" banner row, now
# becomes a normal data row with a new label + real numbers
# ---------------------------------------------------------------------------
$ws.Range("A5:H5").UnMerge()
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 1226
$ws.Range("C5").Value = 1221
$ws.Range("D5").Value = 1201
$ws.Range("E5").Value = 1227
$ws.Range("F5").Value = 1239
$ws.Range("G5").Value = 1231
$ws.Range("H5").Value = 1225
$ws.Range("I5").Value = 1204
$ws.Range("B5:I5").NumberFormat = "#\ ##0"
$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------------------
# Row 6: new row holding the "Source:" note (previously on row 5) and the old
# confidential-data note text is dropped entirely.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = $ws.Range("A6").Value
$ws.Rows.Item(6).RowHeight = 27.75
